# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New strikeout ("K") values computed/scraped for each game row (2..31),
# replacing the previous "Strike#" values that used to live in column G.
$kValues = @{
    2  = 0
    3  = 2
    4  = 3
    5  = 7
    6  = 3
    7  = 3
    8  = 5
    9  = 6
    10 = 2
    11 = 3
    12 = 6
    13 = 5
    14 = 5
    15 = 5
    16 = 5
    17 = 4
    18 = 8
    19 = 3
    20 = 3
    21 = 5
    22 = 3
    23 = 4
    24 = 3
    25 = 4
    26 = 4
    27 = 2
    28 = 3
    29 = 2
    30 = 3
    31 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
